$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Sprint 6 "actual" burn-down data: fill in the actual-points-burned
# table (J25:O.. ) that records which tasks/people burned points on
# which day, then roll the Total: row down to make room and refresh
# the dependent totals/formulas/chart range.
# ------------------------------------------------------------------

# First, make room: the old "Total:" row living at row 39 needs to
# move down to row 40 so the new task row can occupy row 39. Inserting
# a row at 39 shifts the existing row 39 (and its formula) down to 40.
$ws.Rows("39").Insert()

# New task/points-by-person-day rows (row 25 through the row that was
# already holding Q27's blank, left-aligned style cell).
$ws.Range("J25").Value = "18a"
$ws.Range("M25").Value = 1

$ws.Range("J26").Value = "21c"
$ws.Range("M26").Value = 1

$ws.Range("J27").Value = "38a"
$ws.Range("K27").Value = 1

$ws.Range("J28").Value = "38b"
$ws.Range("K28").Value = 1

$ws.Range("J29").Value = 43
$ws.Range("M29").Value = 3

$ws.Range("J30").Value = 46
$ws.Range("M30").Value = 3

$ws.Range("J31").Value = 47
$ws.Range("N31").Value = 5

$ws.Range("J32").Value = 50
$ws.Range("N32").Value = 5

$ws.Range("J33").Value = "52a"
$ws.Range("O33").Value = 5

$ws.Range("J34").Value = "52b"
$ws.Range("N34").Value = 2

$ws.Range("J35").Value = "52c"
$ws.Range("N35").Value = 2

$ws.Range("J36").Value = "52d"
$ws.Range("N36").Value = 2

$ws.Range("J37").Value = "52e"
$ws.Range("N37").Value = 1

$ws.Range("J38").Value = 54
$ws.Range("N38").Value = 1

$ws.Range("J39").Value = "54a"
$ws.Range("N39").Value = 1

# Roll the grand-total formula in (now) row 40 forward to cover the
# newly added row 39.
$ws.Range("K40").Formula = "=SUM(K25:O39)"

# The "Actual" burn row (row 15) sums what got burned each day from
# the new table; extend the ranges summed from row 38 to row 39.
$ws.Range("D15").Formula = "=C15-SUM(K25:K39)"
$ws.Range("E15").Formula = "=D15-SUM(L25:L39)"
$ws.Range("F15").Formula = "=E15-SUM(M25:M39)"
$ws.Range("G15").Formula = "=F15-SUM(N25:N39)"
$ws.Range("H15").Formula = "=G15-SUM(O25:O39)"

# Match the author's final selection before saving.
[void]$ws.Range("P20").Select()

# Update the burn-down chart's value axis to fit the now-negative
# "Actual" series (min/max widen from 0/33 to -5/35).
$chart = $ws.ChartObjects(1).Chart
$chart.Axes(2).MinimumScale = -5
$chart.Axes(2).MaximumScale = 35
